$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("princethorns")

# Un-hide columns A and B and give them explicit widths (previously hidden,
# zero-width helper columns).
$ws.Columns.Item(1).Hidden = $false
$ws.Columns.Item(2).Hidden = $false
$ws.Columns.Item(1).ColumnWidth = 7.33
$ws.Columns.Item(2).ColumnWidth = 10.5

# Make "princethorns" the active/selected sheet (it was "Formatted" before).
$ws.Activate()

Write-Host "done"
